# The workbook tracks a single config sheet. The author:
#  1. Renamed the example STL part file referenced in B4 from "nist.stl" to
#     "nut.stl" (this also causes the now-unused "nist.stl" shared string to
#     be dropped and "nut.stl" appended as a new shared string on save).
#  2. Left the selection/scroll position on B5 instead of B29 (so the view
#     no longer needs to be scrolled down to row 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the part file name value.
$ws.Range("B4").Value = "nut.stl"

# 2. Move the selection/active cell back up to B5 and reset the scrolled
#    top-left cell (selecting a cell that is already on-screen from A1
#    clears any previous "topLeftCell" scroll offset).
$ws.Range("A1").Select()
$ws.Range("B5").Select()
